# Implement image processing memory leaks
# Appends one new trailing data row to each of the four worksheets,
# matching the shape of the existing rows (columns A-I).

$wb = $excel.ActiveWorkbook

# --- Sheet "ROW50-FE-LIFTER": append row 99 ---------------------------------
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$r = 99
$ws1.Cells.Item($r, 1).Value = [double]"45773.81649944445"
$ws1.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item($r, 2).Value = "0x01,0x90"
$ws1.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Cells.Item($r, 4).Value = "0x01,0x36"
$ws1.Cells.Item($r, 5).Value = "0xe"
$ws1.Cells.Item($r, 6).Value = 400
$ws1.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws1.Cells.Item($r, 8).Value = 310
$ws1.Cells.Item($r, 9).Value = 14

# --- Sheet "ROW50-MID-LIFTER": append row 101 -------------------------------
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$r = 101
$ws2.Cells.Item($r, 1).Value = [double]"45773.7771875"
$ws2.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item($r, 2).Value = "0x01,0x90 "
$ws2.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item($r, 4).Value = "0x01,0x3a"
$ws2.Cells.Item($r, 5).Value = "0x19"
$ws2.Cells.Item($r, 6).Value = 400
# Stored as text in the source file (number too large to round-trip as a
# double without losing the exact digit string), so force text entry with a
# leading apostrophe and then strip the resulting quote-prefix style back to
# the default so no stray style index is introduced.
$ws2.Cells.Item($r, 7).Value = "'568631262647113771663628"
$ws2.Cells.Item($r, 7).Style = "Normal"
$ws2.Cells.Item($r, 8).Value = 314
$ws2.Cells.Item($r, 9).Value = 25

# --- Sheet "ROW11-FE-LIFTER": append row 99 ---------------------------------
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$r = 99
$ws3.Cells.Item($r, 1).Value = [double]"45773.84613034722"
$ws3.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item($r, 2).Value = "0x01,0x90"
$ws3.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Cells.Item($r, 4).Value = "0x01,0x36"
$ws3.Cells.Item($r, 5).Value = "0x14"
$ws3.Cells.Item($r, 6).Value = 400
$ws3.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws3.Cells.Item($r, 8).Value = 310
$ws3.Cells.Item($r, 9).Value = 20

# --- Sheet "ROW11-MID-LIFTER": append row 99 --------------------------------
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$r = 99
$ws4.Cells.Item($r, 1).Value = [double]"45773.96371234953"
$ws4.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Cells.Item($r, 2).Value = "0x01,0x90"
$ws4.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Cells.Item($r, 4).Value = "0x01,0x3e"
$ws4.Cells.Item($r, 5).Value = "0x19"
$ws4.Cells.Item($r, 6).Value = 400
$ws4.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws4.Cells.Item($r, 8).Value = 318
$ws4.Cells.Item($r, 9).Value = 25
